# feat: add 2022-Q4 data
#
# Plan:
#  1. "总计" (summary) sheet: insert a new "2022-Q4" entry as the new first
#     data row, push the existing "2022-Q3" / "2022-Q1" rows down by one.
#  2. Insert a brand-new worksheet named "2022-Q4" right before the existing
#     "2022-Q3" sheet, and populate it with the per-fund holdings table.
#  3. Restore the originally-active sheet/tab ("2022-Q1") so the view state
#     matches what it was before the edit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" summary sheet.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Remember the current (pre-edit) row 2 / row 3 contents before overwriting.
$q3Label = $summary.Cells.Item(2, 2).Value2
$q3Count = $summary.Cells.Item(2, 3).Value2
$q3Value = $summary.Cells.Item(2, 4).Value2

$q1Label = $summary.Cells.Item(3, 2).Value2
$q1Count = $summary.Cells.Item(3, 3).Value2
$q1Value = $summary.Cells.Item(3, 4).Value2

# Row 4 becomes the old row 3 ("2022-Q1"), index 2.
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = $q1Label
$summary.Cells.Item(4, 3).Value = $q1Count
$summary.Cells.Item(4, 4).Value = $q1Value

# Give the new A4 cell the same look (border/font) as the other index cells.
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(4, 1).PasteSpecial(-4122)   # xlPasteFormats
$summary.Cells.Item(4, 1).Value = 2
$excel.CutCopyMode = 0

# Row 3 becomes the old row 2 ("2022-Q3"), index 1 (unchanged index value).
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = $q3Label
$summary.Cells.Item(3, 3).Value = $q3Count
$summary.Cells.Item(3, 4).Value = $q3Value

# Row 2 becomes the brand new "2022-Q4" entry, index 0 (unchanged index value).
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 4
$summary.Cells.Item(2, 4).Value = 0.59

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q4" worksheet before "2022-Q3".
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3Sheet)
$q4.Name = "2022-Q4"

$q4.PageSetup.LeftMargin = 0.75 * 72
$q4.PageSetup.RightMargin = 0.75 * 72
$q4.PageSetup.TopMargin = 1 * 72
$q4.PageSetup.BottomMargin = 1 * 72
$q4.PageSetup.HeaderMargin = 0.5 * 72
$q4.PageSetup.FooterMargin = 0.5 * 72

# Header row.
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

# Match the header formatting used on the sibling sheets (bordered/bold/centered style).
$q3Sheet.Cells.Item(1, 2).Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

# Index column (A) style, matching the sibling sheets.
$q3Sheet.Cells.Item(2, 1).Copy()
$q4.Range("A2:A5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Force text storage for the numeric-looking columns. D:F are text on every
# data row; G is text except where the underlying value is a bare zero (same
# convention used on the sibling "2022-Q3"/"2022-Q1" sheets).
$q4.Range("B2:F5").NumberFormat = "@"
$q4.Range("G2:G3").NumberFormat = "@"

# Row 2 - 华安安华灵活配置混合A
$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "002350"
$q4.Cells.Item(2, 3).Value = "华安安华灵活配置混合A"
$q4.Cells.Item(2, 4).Value = "28.67"
$q4.Cells.Item(2, 5).Value = "94.34"
$q4.Cells.Item(2, 6).Value = "1.99"
$q4.Cells.Item(2, 7).Value = "0.5705"
$q4.Cells.Item(2, 8).Value = 7

# Row 3 - 华安安华灵活配置混合C
$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Value = "016183"
$q4.Cells.Item(3, 3).Value = "华安安华灵活配置混合C"
$q4.Cells.Item(3, 4).Value = "0.98"
$q4.Cells.Item(3, 5).Value = "94.34"
$q4.Cells.Item(3, 6).Value = "1.99"
$q4.Cells.Item(3, 7).Value = "0.0195"
$q4.Cells.Item(3, 8).Value = 7

# Row 4 - 恒生前海兴享混合C
$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 2).Value = "014745"
$q4.Cells.Item(4, 3).Value = "恒生前海兴享混合C"
$q4.Cells.Item(4, 4).Value = "0.00"
$q4.Cells.Item(4, 5).Value = "82.60"
$q4.Cells.Item(4, 6).Value = "5.03"
$q4.Cells.Item(4, 7).Value = 0
$q4.Cells.Item(4, 8).Value = 7

# Row 5 - 恒生前海兴享混合A
$q4.Cells.Item(5, 1).Value = 3
$q4.Cells.Item(5, 2).Value = "014744"
$q4.Cells.Item(5, 3).Value = "恒生前海兴享混合A"
$q4.Cells.Item(5, 4).Value = "0.00"
$q4.Cells.Item(5, 5).Value = "82.60"
$q4.Cells.Item(5, 6).Value = "5.03"
$q4.Cells.Item(5, 7).Value = 0
$q4.Cells.Item(5, 8).Value = 7

# ---------------------------------------------------------------------
# 3. Restore the originally-active sheet/tab.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
